$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 12.55417378794944
$ws.Cells.Item(2, 3).Value = 7.107745258394053
$ws.Cells.Item(2, 4).Value = 4.954949426048415
$ws.Cells.Item(2, 5).Value = 7.375301174842367
$ws.Cells.Item(2, 6).Value = 25.14060664150661
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 3.793716826257285
$ws.Cells.Item(2, 9).Value = 4.313184277558927
$ws.Cells.Item(2, 10).Value = 4.251792669778333
$ws.Cells.Item(2, 11).Value = 18.52929954462233
$ws.Cells.Item(2, 12).Value = 5.993327751206099
$ws.Cells.Item(2, 13).Value = 10.73305992979829
$ws.Cells.Item(2, 14).Value = 6.324085418560286
$ws.Cells.Item(2, 15).Value = 9.548466807529088
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 19.10271711392604
$ws.Cells.Item(3, 2).Value = 11.73387998237361
$ws.Cells.Item(3, 3).Value = 6.876544480545292
$ws.Cells.Item(3, 4).Value = 4.783273884031394
$ws.Cells.Item(3, 5).Value = 7.274908603713707
$ws.Cells.Item(3, 6).Value = 24.99396959955318
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 3.982551020108033
$ws.Cells.Item(3, 9).Value = 4.463523741918014
$ws.Cells.Item(3, 10).Value = 4.251792669778333
$ws.Cells.Item(3, 11).Value = 18.59558688050084
$ws.Cells.Item(3, 12).Value = 5.96178347372646
$ws.Cells.Item(3, 13).Value = 10.14041427821327
$ws.Cells.Item(3, 14).Value = 6.193385441098139
$ws.Cells.Item(3, 15).Value = 9.19223243323809
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 19.11055675299681
$ws.Cells.Item(4, 2).Value = 11.19863523783709
$ws.Cells.Item(4, 3).Value = 6.731347330218459
$ws.Cells.Item(4, 4).Value = 4.675760252695765
$ws.Cells.Item(4, 5).Value = 7.211811077695258
$ws.Cells.Item(4, 6).Value = 24.91009986908023
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 4.102838123712686
$ws.Cells.Item(4, 9).Value = 4.559896694139598
$ws.Cells.Item(4, 10).Value = 4.251792669778333
$ws.Cells.Item(4, 11).Value = 18.63952607937159
$ws.Cells.Item(4, 12).Value = 5.941488826190406
$ws.Cells.Item(4, 13).Value = 9.758049892784443
$ws.Cells.Item(4, 14).Value = 6.112529093527351
$ws.Cells.Item(4, 15).Value = 8.96723759143487
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 19.12042970713275
$ws.Cells.Item(5, 2).Value = 10.97655228318597
$ws.Cells.Item(5, 3).Value = 6.675846332497862
$ws.Cells.Item(5, 4).Value = 4.632572711992714
$ws.Cells.Item(5, 5).Value = 7.184872934079014
$ws.Cells.Item(5, 6).Value = 24.86940736779328
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 4.153275638880205
$ws.Cells.Item(5, 9).Value = 4.602501158099312
$ws.Cells.Item(5, 10).Value = 4.251792669778333
$ws.Cells.Item(5, 11).Value = 18.65250306892706
$ws.Cells.Item(5, 12).Value = 5.932546698074443
$ws.Cells.Item(5, 13).Value = 9.599097954492974
$ws.Cells.Item(5, 14).Value = 6.080711387048089
$ws.Cells.Item(5, 15).Value = 8.874523884662183
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 19.11992186391201
$ws.Cells.Item(6, 2).Value = 10.93716052176938
$ws.Cells.Item(6, 3).Value = 6.672065265874346
$ws.Cells.Item(6, 4).Value = 4.62671849511989
$ws.Cells.Item(6, 5).Value = 7.179317592399487
$ws.Cells.Item(6, 6).Value = 24.85288976684437
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 4.16210121534563
$ws.Cells.Item(6, 9).Value = 4.612569941680051
$ws.Cells.Item(6, 10).Value = 4.251792669778333
$ws.Cells.Item(6, 11).Value = 18.64770688090177
$ws.Cells.Item(6, 12).Value = 5.930514772000388
$ws.Cells.Item(6, 13).Value = 9.574196448054119
$ws.Cells.Item(6, 14).Value = 6.076961067319625
$ws.Cells.Item(6, 15).Value = 8.859577087866565
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 19.11285174996479
$ws.Cells.Item(7, 2).Value = 11.18858667659944
$ws.Cells.Item(7, 3).Value = 6.745448520504195
$ws.Cells.Item(7, 4).Value = 4.678803567848899
$ws.Cells.Item(7, 5).Value = 7.208578845651652
$ws.Cells.Item(7, 6).Value = 24.88274060500341
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 4.104541044130533
$ws.Cells.Item(7, 9).Value = 4.56821784834304
$ws.Cells.Item(7, 10).Value = 4.251792669778333
$ws.Cells.Item(7, 11).Value = 18.62062320953746
$ws.Cells.Item(7, 12).Value = 5.939934731007279
$ws.Cells.Item(7, 13).Value = 9.760667963348348
$ws.Cells.Item(7, 14).Value = 6.116275844712284
$ws.Cells.Item(7, 15).Value = 8.967421142237377
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 19.10118779169014
$ws.Cells.Item(8, 2).Value = 12.26932579603854
$ws.Cells.Item(8, 3).Value = 7.047842521217975
$ws.Cells.Item(8, 4).Value = 4.900811965224364
$ws.Cells.Item(8, 5).Value = 7.337315597640674
$ws.Cells.Item(8, 6).Value = 25.05382367452751
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 3.859207679878188
$ws.Cells.Item(8, 9).Value = 4.373953509657684
$ws.Cells.Item(8, 10).Value = 4.251792669778333
$ws.Cells.Item(8, 11).Value = 18.52644433553998
$ws.Cells.Item(8, 12).Value = 5.980834289606944
$ws.Cells.Item(8, 13).Value = 10.53850665756945
$ws.Cells.Item(8, 14).Value = 6.284554585893996
$ws.Cells.Item(8, 15).Value = 9.42877055375083
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 19.07916430862167
$ws.Cells.Item(9, 2).Value = 14.1503979539784
$ws.Cells.Item(9, 3).Value = 7.588820660330501
$ws.Cells.Item(9, 4).Value = 5.311171817954766
$ws.Cells.Item(9, 5).Value = 7.582780760298931
$ws.Cells.Item(9, 6).Value = 25.48553143729028
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 3.410260842059372
$ws.Cells.Item(9, 9).Value = 4.012174383039591
$ws.Cells.Item(9, 10).Value = 4.251792669778333
$ws.Cells.Item(9, 11).Value = 18.40188977867361
$ws.Cells.Item(9, 12).Value = 6.056196896649164
$ws.Cells.Item(9, 13).Value = 11.90724297847367
$ws.Cells.Item(9, 14).Value = 6.601362365684847
$ws.Cells.Item(9, 15).Value = 10.27709685807912
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 19.11136961751946
$ws.Cells.Item(10, 2).Value = 15.37589204112835
$ws.Cells.Item(10, 3).Value = 8.000628033423405
$ws.Cells.Item(10, 4).Value = 5.57219997998118
$ws.Cells.Item(10, 5).Value = 7.692686691058285
$ws.Cells.Item(10, 6).Value = 25.68135194812274
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 3.124414113841304
$ws.Cells.Item(10, 9).Value = 3.774334580297831
$ws.Cells.Item(10, 10).Value = 4.251792669778333
$ws.Cells.Item(10, 11).Value = 18.24243910750795
$ws.Cells.Item(10, 12).Value = 6.094458408414714
$ws.Cells.Item(10, 13).Value = 12.83304867705525
$ws.Cells.Item(10, 14).Value = 6.771629503265272
$ws.Cells.Item(10, 15).Value = 10.81856453860499
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 19.06648657162475
$ws.Cells.Item(11, 2).Value = 15.86384845949349
$ws.Cells.Item(11, 3).Value = 8.430239842875332
$ws.Cells.Item(11, 4).Value = 5.451174545603521
$ws.Cells.Item(11, 5).Value = 7.264448956287217
$ws.Cells.Item(11, 6).Value = 24.54191984529939
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 3.775815705518889
$ws.Cells.Item(11, 9).Value = 3.722108202015785
$ws.Cells.Item(11, 10).Value = 4.251792669778333
$ws.Cells.Item(11, 11).Value = 17.5045184599554
$ws.Cells.Item(11, 12).Value = 6.046490227149508
$ws.Cells.Item(11, 13).Value = 13.33568387309926
$ws.Cells.Item(11, 14).Value = 6.35564569259995
$ws.Cells.Item(11, 15).Value = 10.67303383675447
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 18.28200469801562
$ws.Cells.Item(12, 2).Value = 16.03604160826671
$ws.Cells.Item(12, 3).Value = 8.702556389532006
$ws.Cells.Item(12, 4).Value = 5.286612884295645
$ws.Cells.Item(12, 5).Value = 6.969594885553615
$ws.Cells.Item(12, 6).Value = 23.57787469441139
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 4.88669651886272
$ws.Cells.Item(12, 9).Value = 3.710273992870889
$ws.Cells.Item(12, 10).Value = 4.251792669778333
$ws.Cells.Item(12, 11).Value = 16.96263403504058
$ws.Cells.Item(12, 12).Value = 6.056712772415045
$ws.Cells.Item(12, 13).Value = 13.56810398284468
$ws.Cells.Item(12, 14).Value = 5.982045659979262
$ws.Cells.Item(12, 15).Value = 10.42960571119377
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 17.66219964572736
$ws.Cells.Item(13, 2).Value = 15.97355858783464
$ws.Cells.Item(13, 3).Value = 8.893343664654447
$ws.Cells.Item(13, 4).Value = 5.080119672738717
$ws.Cells.Item(13, 5).Value = 6.763843733871646
$ws.Cells.Item(13, 6).Value = 22.64980208793832
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 6.154989831156969
$ws.Cells.Item(13, 9).Value = 3.741388769410038
$ws.Cells.Item(13, 10).Value = 4.251792669778333
$ws.Cells.Item(13, 11).Value = 16.51382657442097
$ws.Cells.Item(13, 12).Value = 6.107478106046078
$ws.Cells.Item(13, 13).Value = 13.61904682156461
$ws.Cells.Item(13, 14).Value = 5.62617721591972
$ws.Cells.Item(13, 15).Value = 10.09291275195208
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 17.10434068037842
$ws.Cells.Item(14, 2).Value = 15.82129687118618
$ws.Cells.Item(14, 3).Value = 8.997359905467221
$ws.Cells.Item(14, 4).Value = 4.915853825694412
$ws.Cells.Item(14, 5).Value = 6.67735761894312
$ws.Cells.Item(14, 6).Value = 22.01637654197448
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 7.092128591422657
$ws.Cells.Item(14, 9).Value = 3.784972165313909
$ws.Cells.Item(14, 10).Value = 4.251792669778333
$ws.Cells.Item(14, 11).Value = 16.24513493959228
$ws.Cells.Item(14, 12).Value = 6.166968342291035
$ws.Cells.Item(14, 13).Value = 13.57457809761914
$ws.Cells.Item(14, 14).Value = 5.392935300455621
$ws.Cells.Item(14, 15).Value = 9.812138910586409
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 16.74293878483909
$ws.Cells.Item(15, 2).Value = 15.73113612995534
$ws.Cells.Item(15, 3).Value = 9.01016961847386
$ws.Cells.Item(15, 4).Value = 4.868170985620058
$ws.Cells.Item(15, 5).Value = 6.663040175553276
$ws.Cells.Item(15, 6).Value = 21.86269316126042
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 7.315353441610852
$ws.Cells.Item(15, 9).Value = 3.808858431369081
$ws.Cells.Item(15, 10).Value = 4.251792669778333
$ws.Cells.Item(15, 11).Value = 16.19275721007118
$ws.Cells.Item(15, 12).Value = 6.182465711703395
$ws.Cells.Item(15, 13).Value = 13.52600002048779
$ws.Cells.Item(15, 14).Value = 5.338649659207335
$ws.Cells.Item(15, 15).Value = 9.724516529331799
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 16.66263627113459
$ws.Cells.Item(16, 2).Value = 15.23726144946548
$ws.Cells.Item(16, 3).Value = 8.815647274679593
$ws.Cells.Item(16, 4).Value = 4.793500840214279
$ws.Cells.Item(16, 5).Value = 6.651939229967577
$ws.Cells.Item(16, 6).Value = 21.94773023862912
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 7.153865862722808
$ws.Cells.Item(16, 9).Value = 3.909196221010635
$ws.Cells.Item(16, 10).Value = 4.251792669778333
$ws.Cells.Item(16, 11).Value = 16.3463873271953
$ws.Cells.Item(16, 12).Value = 6.150244018900806
$ws.Cells.Item(16, 13).Value = 13.14301058638223
$ws.Cells.Item(16, 14).Value = 5.335157724174773
$ws.Cells.Item(16, 15).Value = 9.55194412580542
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 16.78600903912118
$ws.Cells.Item(17, 2).Value = 14.93434213054631
$ws.Cells.Item(17, 3).Value = 8.606640345555975
$ws.Cells.Item(17, 4).Value = 4.825893530288488
$ws.Cells.Item(17, 5).Value = 6.679598611408884
$ws.Cells.Item(17, 6).Value = 22.34691829218641
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 6.508207932883679
$ws.Cells.Item(17, 9).Value = 3.963674641184693
$ws.Cells.Item(17, 10).Value = 4.251792669778333
$ws.Cells.Item(17, 11).Value = 16.60193305889215
$ws.Cells.Item(17, 12).Value = 6.089862617220501
$ws.Cells.Item(17, 13).Value = 12.86775721333897
$ws.Cells.Item(17, 14).Value = 5.453632250488045
$ws.Cells.Item(17, 15).Value = 9.572547005463413
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 17.06941907663512
$ws.Cells.Item(18, 2).Value = 14.77402591999843
$ws.Cells.Item(18, 3).Value = 8.354793030296255
$ws.Cells.Item(18, 4).Value = 4.951423272112622
$ws.Cells.Item(18, 5).Value = 6.799229061670211
$ws.Cells.Item(18, 6).Value = 23.0889322083519
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 5.430593762480435
$ws.Cells.Item(18, 9).Value = 3.974490215775643
$ws.Cells.Item(18, 10).Value = 4.251792669778333
$ws.Cells.Item(18, 11).Value = 17.00137127391113
$ws.Cells.Item(18, 12).Value = 6.024701977646707
$ws.Cells.Item(18, 13).Value = 12.65638263207773
$ws.Cells.Item(18, 14).Value = 5.706635764185116
$ws.Cells.Item(18, 15).Value = 9.760340494605801
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 17.54355123560419
$ws.Cells.Item(19, 2).Value = 14.73308741926411
$ws.Cells.Item(19, 3).Value = 8.129671130030191
$ws.Cells.Item(19, 4).Value = 5.146536596951663
$ws.Cells.Item(19, 5).Value = 7.057007390006505
$ws.Cells.Item(19, 6).Value = 24.02651752432506
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 4.250762681927696
$ws.Cells.Item(19, 9).Value = 3.962812944478
$ws.Cells.Item(19, 10).Value = 4.251792669778333
$ws.Cells.Item(19, 11).Value = 17.48363299893108
$ws.Cells.Item(19, 12).Value = 6.000507586596223
$ws.Cells.Item(19, 13).Value = 12.52405476893005
$ws.Cells.Item(19, 14).Value = 6.0782916289296
$ws.Cells.Item(19, 15).Value = 10.06676110696505
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(19, 17).Value = 18.12065743503108
$ws.Cells.Item(20, 2).Value = 15.04711345738298
$ws.Cells.Item(20, 3).Value = 7.937281989874642
$ws.Cells.Item(20, 4).Value = 5.512966123779141
$ws.Cells.Item(20, 5).Value = 7.653546128937235
$ws.Cells.Item(20, 6).Value = 25.54301271899249
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 3.200786115222866
$ws.Cells.Item(20, 9).Value = 3.860899671457827
$ws.Cells.Item(20, 10).Value = 4.251792669778333
$ws.Cells.Item(20, 11).Value = 18.22239091211992
$ws.Cells.Item(20, 12).Value = 6.080360433317381
$ws.Cells.Item(20, 13).Value = 12.60898223546369
$ws.Cells.Item(20, 14).Value = 6.736349824839048
$ws.Cells.Item(20, 15).Value = 10.68046228893921
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(20, 17).Value = 19.01504125974067
$ws.Cells.Item(21, 2).Value = 15.95054260398934
$ws.Cells.Item(21, 3).Value = 8.199413697918766
$ws.Cells.Item(21, 4).Value = 5.749112254639786
$ws.Cells.Item(21, 5).Value = 7.832581144658814
$ws.Cells.Item(21, 6).Value = 25.93631485467559
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 2.950837243936141
$ws.Cells.Item(21, 9).Value = 3.671943356869769
$ws.Cells.Item(21, 10).Value = 4.251792669778333
$ws.Cells.Item(21, 11).Value = 18.23565238346877
$ws.Cells.Item(21, 12).Value = 6.128313320735635
$ws.Cells.Item(21, 13).Value = 13.26943927400087
$ws.Cells.Item(21, 14).Value = 6.954300227129065
$ws.Cells.Item(21, 15).Value = 11.15090621108229
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = 19.13634285271894
$ws.Cells.Item(22, 2).Value = 16.51871463373844
$ws.Cells.Item(22, 3).Value = 8.363646360925232
$ws.Cells.Item(22, 4).Value = 5.882207233338763
$ws.Cells.Item(22, 5).Value = 7.919815824803732
$ws.Cells.Item(22, 6).Value = 26.15971049381766
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 2.801299046818165
$ws.Cells.Item(22, 9).Value = 3.543744850856595
$ws.Cells.Item(22, 10).Value = 4.251792669778333
$ws.Cells.Item(22, 11).Value = 18.23301881094863
$ws.Cells.Item(22, 12).Value = 6.154189038913568
$ws.Cells.Item(22, 13).Value = 13.6850267752882
$ws.Cells.Item(22, 14).Value = 7.060235348531485
$ws.Cells.Item(22, 15).Value = 11.42619478157142
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = 19.20064772894596
$ws.Cells.Item(23, 2).Value = 16.22389245840881
$ws.Cells.Item(23, 3).Value = 8.260830049214222
$ws.Cells.Item(23, 4).Value = 5.807977917235991
$ws.Cells.Item(23, 5).Value = 7.876342007493728
$ws.Cells.Item(23, 6).Value = 26.06923895825595
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 2.879881404014507
$ws.Cells.Item(23, 9).Value = 3.602020579102235
$ws.Cells.Item(23, 10).Value = 4.251792669778333
$ws.Cells.Item(23, 11).Value = 18.25580204016098
$ws.Cells.Item(23, 12).Value = 6.141834564543773
$ws.Cells.Item(23, 13).Value = 13.46044500199159
$ws.Cells.Item(23, 14).Value = 6.999437418612978
$ws.Cells.Item(23, 15).Value = 11.27875051119036
$ws.Cells.Item(23, 16).Value = 0
$ws.Cells.Item(23, 17).Value = 19.18740743772729
$ws.Cells.Item(24, 2).Value = 15.04450127172835
$ws.Cells.Item(24, 3).Value = 7.888628946396816
$ws.Cells.Item(24, 4).Value = 5.525810135875967
$ws.Cells.Item(24, 5).Value = 7.704179294639564
$ws.Cells.Item(24, 6).Value = 25.69281831667354
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 3.18553380813438
$ws.Cells.Item(24, 9).Value = 3.844884872709858
$ws.Cells.Item(24, 10).Value = 4.251792669778333
$ws.Cells.Item(24, 11).Value = 18.31560323218736
$ws.Cells.Item(24, 12).Value = 6.091476815855118
$ws.Cells.Item(24, 13).Value = 12.58205107296206
$ws.Cells.Item(24, 14).Value = 6.773254552273772
$ws.Cells.Item(24, 15).Value = 10.70618082721889
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 19.11510112979267
$ws.Cells.Item(25, 2).Value = 13.65749826468535
$ws.Cells.Item(25, 3).Value = 7.470621338303567
$ws.Cells.Item(25, 4).Value = 5.208359553866967
$ws.Cells.Item(25, 5).Value = 7.512883228944984
$ws.Cells.Item(25, 6).Value = 25.31678061989039
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 3.52978163364371
$ws.Cells.Item(25, 9).Value = 4.120484610601053
$ws.Cells.Item(25, 10).Value = 4.251792669778333
$ws.Cells.Item(25, 11).Value = 18.39886901950338
$ws.Cells.Item(25, 12).Value = 6.034249570682294
$ws.Cells.Item(25, 13).Value = 11.56061840499813
$ws.Cells.Item(25, 14).Value = 6.523570243998527
$ws.Cells.Item(25, 15).Value = 10.05615449854039
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(25, 17).Value = 19.06412776056488
